# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.671.81"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$ws.Range("D3").Value = "1.962.06"
$ws.Range("E3").Value = "  +2.23%  "

# Row 4
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "'248.40"
$ws.Range("E5").Value = "  +1.44%  "

# Row 6
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.12%  "

# Row 7
$ws.Range("D7").Value = "0.4819"
$ws.Range("E7").Value = "  -1.00%  "

# Row 8
$ws.Range("E8").Value = "  +0.42%  "

# Row 9
$ws.Range("D9").Value = "'0.2930"
$ws.Range("E9").Value = "  +1.28%  "

# Row 10
$ws.Range("D10").Value = "'0.06750"
$ws.Range("E10").Value = "  +0.66%  "

# Row 11
$ws.Range("D11").Value = "109.19"
$ws.Range("E11").Value = "  -1.96%  "

# Row 12
$ws.Range("D12").Value = "19.18"
$ws.Range("E12").Value = "  +0.94%  "

# Row 13
$ws.Range("D13").Value = "1.959.16"
$ws.Range("E13").Value = "  +2.06%  "

# Row 14
$ws.Range("D14").Value = "0.07736"
$ws.Range("E14").Value = "  +1.96%  "

# Row 15
$ws.Range("D15").Value = "5.456"
$ws.Range("E15").Value = "  +3.36%  "

# Row 16
$ws.Range("D16").Value = "0.6957"
$ws.Range("E16").Value = "  +4.07%  "

# Row 17
$ws.Range("D17").Value = "292.12"
$ws.Range("E17").Value = "  -0.56%  "

# Row 18
$ws.Range("D18").Value = "30.673.28"
$ws.Range("E18").Value = "  +0.42%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "5.662"
$ws.Range("E19").Value = "  +2.06%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "13.16"
$ws.Range("E20").Value = "  +1.23%  "

# Row 21
$ws.Range("D21").Value = "'0.000007690"
$ws.Range("E21").Value = "  +1.61%  "

# Row 22
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").Value = "0.4899"
$ws.Range("E22").Value = "  +12.52%  "

# Row 23
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.209.18"
$ws.Range("E23").Value = "  +1.88%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "0.9987"

# Row 25
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "0.9989"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("B26").Value = "Chainlink"
$ws.Range("C26").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D26").Value = "6.619"
$ws.Range("E26").Value = "  +2.75%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.914"
$ws.Range("E27").Value = "  +4.67%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'169.50"
$ws.Range("E28").Value = "  +2.91%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.98"
$ws.Range("E29").Value = "  -1.26%  "

# Row 30
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "2.175"
$ws.Range("E30").Value = "  +4.20%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1066"
$ws.Range("E31").Value = "  -0.64%  "

# Row 32
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.440"
$ws.Range("E32").Value = "  -0.50%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.766"
$ws.Range("E33").Value = "  +17.58%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.440"
$ws.Range("E34").Value = "  +7.25%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.05091"
$ws.Range("E35").Value = "  +1.47%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7710"
$ws.Range("E36").Value = "  +4.11%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.178"
$ws.Range("E37").Value = "  +3.75%  "

# Row 38
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.727"
$ws.Range("E38").Value = "  +0.47%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02033"
$ws.Range("E39").Value = "  +0.26%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.716"
$ws.Range("E40").Value = "  +1.33%  "

# Row 41
$ws.Range("D41").Value = "6.483"
$ws.Range("E41").Value = "  +11.09%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "2.138"
$ws.Range("E42").Value = "  +6.12%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'110.00"
$ws.Range("E43").Value = "  -0.40%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8848"
$ws.Range("E44").Value = "  +2.14%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4460"
$ws.Range("E45").Value = "  +0.67%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "69.98"
$ws.Range("E46").Value = "  -1.71%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.09%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.464"
$ws.Range("E48").Value = "  +3.35%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1276"
$ws.Range("E49").Value = "  +3.53%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.331"
$ws.Range("E50").Value = "  +1.73%  "

# Row 51
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'36.00"
$ws.Range("E51").Value = "  +3.17%  "

